$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: Target Species "HCON", offset/spatial smooth both "X"
$ws.Range("A27").Value = "HCON"
$ws.Range("B27").Value = "X"
$ws.Range("C27").Value = "X"
$ws.Range("J27").Value = 0.391
$ws.Range("K27").Value = 10.9
$ws.Range("L27").Value = 807.7
$ws.Range("M27").Value = 922.1

# Row 28
$ws.Range("B28").Value = "X"
$ws.Range("C28").Value = "X"
$ws.Range("D28").Value = "X"
$ws.Range("J28").Value = 0.514
$ws.Range("K28").Value = 32
$ws.Range("L28").Value = 450.46
$ws.Range("M28").Value = 885.6

# Row 29 (tall row - wrapped landscape metrics text)
$ws.Range("B29").Value = "X"
$ws.Range("E29").Value = "Agriculture, Forest, Developed"
$ws.Range("F29").Value = "X"
$ws.Range("G29").Value = "X"
$ws.Range("J29").Value = 0.488
$ws.Range("K29").Value = 20.6
$ws.Range("L29").Value = 444.37
$ws.Range("M29").Value = 881.5
$ws.Rows.Item(29).RowHeight = 29

# Row 30 (tall row - wrapped landscape metrics text)
$ws.Range("B30").Value = "X"
$ws.Range("E30").Value = "Agriculture, Forest, Developed"
$ws.Range("F30").Value = "X"
$ws.Range("H30").Value = "X"
$ws.Range("I30").Value = "X"
$ws.Range("J30").Value = 0.308
$ws.Range("K30").Value = 23.1
$ws.Range("L30").Value = 439.3
$ws.Range("M30").Value = 870.9
$ws.Rows.Item(30).RowHeight = 29

# Row 31
$ws.Range("B31").Value = "X"
$ws.Range("E31").Value = "Forest, Developed"
$ws.Range("F31").Value = "X"
$ws.Range("H31").Value = "X"
$ws.Range("I31").Value = "X"
$ws.Range("J31").Value = 0.255
$ws.Range("K31").Value = 22.3
$ws.Range("L31").Value = 440.66
$ws.Range("M31").Value = 870.4

# Row 32
$ws.Range("B32").Value = "X"
$ws.Range("E32").Value = "Agriculture, Forest"
$ws.Range("H32").Value = "X"
$ws.Range("I32").Value = "X"

# Row 33
$ws.Range("B33").Value = "X"
$ws.Range("E33").Value = "Forest"
$ws.Range("H33").Value = "X"
$ws.Range("I33").Value = "X"

# Row 34
$ws.Range("B34").Value = "X"
$ws.Range("H34").Value = "X"
$ws.Range("I34").Value = "X"

# Update the view selection to reflect the new active cell (matches post-edit state)
$ws.Range("L37").Select()
